$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.740.42'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').Value = '2.401.68'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.993'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.62%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.37'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.31'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.63%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').Value = '2.382.53'
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.108'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.161'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('E12').Value = '  -2.27%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.340'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000170'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '2.795.53'
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').Value = '60.697.73'
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = '2.377.80'
$ws.Range('E18').Value = '  -2.33%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.54'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.25'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '321.36'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.95%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.02'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  -6.21%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '64.37'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.13%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.57'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -8.19%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '570.46'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -7.07%  '
$ws.Range('D29').Value = '2.504.23'
$ws.Range('E29').Value = '  -2.43%  '
$ws.Range('D30').Value = '0.0₃0911'
$ws.Range('E30').Value = '  -3.87%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.85'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('E32').Value = '  -5.88%  '
$ws.Range('E33').Value = '  -2.37%  '
$ws.Range('E34').Value = '  -5.64%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.01'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.61'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -5.57%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.368'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.38'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '147.69'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.94%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.15'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.07'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.16%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.66'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.38%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '40.86'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.88%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.33'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.82%  '
$ws.Range('D46').Value = '0.0₆0284'
$ws.Range('E46').Value = '  +19.60%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '140.47'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('E48').Value = '  -3.58%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.585'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.17%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0503'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.77%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.27'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.81%  '
